$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 423
$ws.Range("I6").Value = 423
$ws.Range("K6").Value = 1269
$ws.Range("M6").Value = -1157

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("N26").Value = 0
$ws.Range("L26").ClearContents()
$ws.Range("M26").ClearContents()

$ws.Range("H33").Value = 37953.57
$ws.Range("I33").Value = 50730.7
$ws.Range("J33").Value = 6010.75
$ws.Range("K33").Value = 50730.7
$ws.Range("L33").Value = 6010.75
$ws.Range("M33").Value = -50501.7
$ws.Range("N33").Value = -6468.75

$ws.Range("H40").Value = 6980.4375
$ws.Range("I40").Value = 2912.75
$ws.Range("J40").Value = 11048.125
$ws.Range("K40").Value = 2912.75
$ws.Range("L40").Value = 11048.125
$ws.Range("M40").Value = -2737.75
$ws.Range("N40").Value = -11398.125

$ws.Range("H62").Value = 3479.8
$ws.Range("I62").Value = 3474.75
$ws.Range("K62").Value = 3474.75
$ws.Range("M62").Value = -2850.75

$ws.Range("H65").Value = 3479.8
$ws.Range("I65").Value = 3474.75
$ws.Range("K65").Value = 17373.75
$ws.Range("M65").Value = -14253.75

$ws.Range("H69").Value = 17165.389
$ws.Range("I69").Value = 8999.5
$ws.Range("J69").Value = 18186.125
$ws.Range("K69").Value = 26998.5
$ws.Range("L69").Value = 54558.375
$ws.Range("M69").Value = -26124.5
$ws.Range("N69").Value = -56306.375

$ws.Range("H72").Value = 17165.389
$ws.Range("I72").Value = 8999.5
$ws.Range("J72").Value = 18186.125
$ws.Range("K72").Value = 80995.5
$ws.Range("L72").Value = 163675.125
$ws.Range("M72").Value = -76627.5
$ws.Range("N72").Value = -172411.125

$ws.Range("H98").Value = 368.38095
$ws.Range("I98").Value = 387.70587
$ws.Range("K98").Value = 387.70587
$ws.Range("M98").Value = 1110.29413

$ws.Range("H101").Value = 554.5789
$ws.Range("I101").Value = 629.6667
$ws.Range("K101").Value = 1889.0001
$ws.Range("M101").Value = -267.0001

$ws.Range("H122").Value = 368.38095
$ws.Range("I122").Value = 387.70587
$ws.Range("K122").Value = 1163.11761
$ws.Range("M122").Value = 1286.88239

$ws.Range("H137").Value = 3108.9167
$ws.Range("I137").Value = 2775.1538
$ws.Range("J137").Value = 3503.3635
$ws.Range("K137").Value = 8325.4614
$ws.Range("L137").Value = 10510.0905
$ws.Range("M137").Value = -5775.4614
$ws.Range("N137").Value = -15610.0905

$ws.Range("H141").Value = 4084.2727
$ws.Range("I141").Value = 3013
$ws.Range("K141").Value = 9039
$ws.Range("M141").Value = -3859

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4304.2
$ws.Range("J2").Value = 4499.5
$ws.Range("L2").Value = 4499.5
$ws.Range("N2").Value = -4725.5

$ws.Range("H39").Value = 25303.125
$ws.Range("I39").Value = 10985.714
$ws.Range("J39").Value = 125525
$ws.Range("K39").Value = 10985.714
$ws.Range("L39").Value = 125525
$ws.Range("M39").Value = -10465.714
$ws.Range("N39").Value = -126565

$ws.Range("H116").Value = 4304.2
$ws.Range("J116").Value = 4499.5
$ws.Range("L116").Value = 4499.5
$ws.Range("N116").Value = -9087.5

$ws.Range("H122").Value = 3641.84
$ws.Range("I122").Value = 3792.0952
$ws.Range("J122").Value = 2853
$ws.Range("K122").Value = 11376.2856
$ws.Range("L122").Value = 8559
$ws.Range("M122").Value = -8926.285600000001
$ws.Range("N122").Value = -13459

$ws.Range("H132").Value = 2840.1936
$ws.Range("I132").Value = 1514.8182
$ws.Range("K132").Value = 4544.4546
$ws.Range("M132").Value = -2014.4546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4304.2
$ws.Range("J3").Value = 4499.5
$ws.Range("L3").Value = 4499.5
$ws.Range("N3").Value = -4727.5

$ws.Range("H99").Value = 1465.35
$ws.Range("I99").Value = 1165.2354
$ws.Range("J99").Value = 3166
$ws.Range("K99").Value = 1165.2354
$ws.Range("L99").Value = 3166
$ws.Range("M99").Value = 332.7646
$ws.Range("N99").Value = -6162

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15347.286
$ws.Range("I31").Value = 8264.143
$ws.Range("J31").Value = 22430.428
$ws.Range("K31").Value = 8264.143
$ws.Range("L31").Value = 22430.428
$ws.Range("M31").Value = -7969.143
$ws.Range("N31").Value = -23020.428

$ws.Range("H34").Value = 15347.286
$ws.Range("I34").Value = 8264.143
$ws.Range("J34").Value = 22430.428
$ws.Range("K34").Value = 8264.143
$ws.Range("L34").Value = 22430.428
$ws.Range("M34").Value = -8062.143
$ws.Range("N34").Value = -22834.428

$ws.Range("H132").Value = 2101.476
$ws.Range("I132").Value = 2101.476
$ws.Range("K132").Value = 6304.428
$ws.Range("M132").Value = -3774.428

$ws.Range("H134").Value = 5376.5264
$ws.Range("I134").Value = 2958.3076
$ws.Range("K134").Value = 8874.9228
$ws.Range("M134").Value = -6339.9228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 4222.5
$ws.Range("I14").Value = 4222.5
$ws.Range("K14").Value = 12667.5
$ws.Range("M14").Value = -12494.5

$ws.Range("H62").Value = 9506
$ws.Range("I62").Value = 9506
$ws.Range("K62").Value = 28518
$ws.Range("M62").Value = -27832

$ws.Range("H65").Value = 9506
$ws.Range("I65").Value = 9506
$ws.Range("K65").Value = 85554
$ws.Range("M65").Value = -82122

$ws.Range("H75").Value = 4117.615
$ws.Range("I75").Value = 956
$ws.Range("J75").Value = 4692.4546
$ws.Range("K75").Value = 2868
$ws.Range("L75").Value = 14077.3638
$ws.Range("M75").Value = -1870
$ws.Range("N75").Value = -16073.3638

$ws.Range("H78").Value = 4117.615
$ws.Range("I78").Value = 956
$ws.Range("J78").Value = 4692.4546
$ws.Range("K78").Value = 8604
$ws.Range("L78").Value = 42232.0914
$ws.Range("M78").Value = -3612
$ws.Range("N78").Value = -52216.0914

$ws.Range("H103").Value = 159.4
$ws.Range("I103").Value = 102.666664
$ws.Range("J103").Value = 244.5
$ws.Range("K103").Value = 307.999992
$ws.Range("L103").Value = 733.5
$ws.Range("M103").Value = 571.000008
$ws.Range("N103").Value = -2491.5

$ws.Range("H107").Value = 350
$ws.Range("J107").Value = 350
$ws.Range("L107").Value = 1050
$ws.Range("N107").Value = -4890

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 16300
$ws.Range("J19").Value = 9900
$ws.Range("L19").Value = 9900
$ws.Range("N19").Value = -10476

$ws.Range("H44").Value = 39999
$ws.Range("J44").Value = 39999
$ws.Range("L44").Value = 39999
$ws.Range("N44").Value = -41191

$ws.Range("H70").Value = 4993.6665
$ws.Range("I70").Value = 4437.25
$ws.Range("K70").Value = 4437.25
$ws.Range("M70").Value = -4167.25

$ws.Range("H73").Value = 4993.6665
$ws.Range("I73").Value = 4437.25
$ws.Range("K73").Value = 4437.25
$ws.Range("M73").Value = -3501.25

$ws.Range("H102").Value = 2704.0293
$ws.Range("I102").Value = 1347.174
$ws.Range("K102").Value = 1347.174
$ws.Range("M102").Value = 274.826

$ws.Range("H113").Value = 3090.3157
$ws.Range("I113").Value = 3278.182
$ws.Range("K113").Value = 3278.182
$ws.Range("M113").Value = -1108.182

$ws.Range("H122").Value = 5132.3184
$ws.Range("I122").Value = 5093.231
$ws.Range("K122").Value = 15279.693
$ws.Range("M122").Value = -12829.693

$ws.Range("H126").Value = 4058
$ws.Range("I126").Value = 2824.8
$ws.Range("J126").Value = 5599.5
$ws.Range("K126").Value = 8474.400000000001
$ws.Range("L126").Value = 16798.5
$ws.Range("M126").Value = -6004.400000000001
$ws.Range("N126").Value = -21738.5

$ws.Range("H138").Value = 94996.75
$ws.Range("J138").Value = 94996.75
$ws.Range("L138").Value = 94996.75
$ws.Range("N138").Value = -105276.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 65079.35
$ws.Range("I7").Value = 72723.266
$ws.Range("K7").Value = 72723.266
$ws.Range("M7").Value = -72611.266

$ws.Range("H40").Value = 3076.2058
$ws.Range("I40").Value = 2389.1785
$ws.Range("K40").Value = 2389.1785
$ws.Range("M40").Value = -2253.1785

$ws.Range("H46").Value = 4523.6665
$ws.Range("J46").Value = 4726.7144
$ws.Range("L46").Value = 4726.7144
$ws.Range("N46").Value = -5102.7144

$ws.Range("H122").Value = 5583.8823
$ws.Range("I122").Value = 4531.6665
$ws.Range("J122").Value = 8109.2
$ws.Range("K122").Value = 13594.9995
$ws.Range("L122").Value = 24327.6
$ws.Range("M122").Value = -11144.9995
$ws.Range("N122").Value = -29227.6

$ws.Range("H126").Value = 65079.35
$ws.Range("I126").Value = 72723.266
$ws.Range("K126").Value = 218169.798
$ws.Range("M126").Value = -215699.798

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4552.933
$ws.Range("I122").Value = 3391.5833
$ws.Range("M122").Value = -7724.749899999999

$ws.Range("H126").Value = 1670.9584
$ws.Range("I126").Value = 1599.3158
$ws.Range("K126").Value = 4797.9474
$ws.Range("M126").Value = -2327.9474
